$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.454345666666666
$ws.Range("H2").Value = 7.363036999999999
$ws.Range("I2").Value = 0.2857469401257222
$ws.Range("J2").Value = 0.3027613480760606
$ws.Range("M2").Value = 36.48539666666667
$ws.Range("N2").Value = 109.45619
$ws.Range("O2").Value = 0.4260639713374229
$ws.Range("P2").Value = 0.4324607845540777
$ws.Range("Q2").Value = 89.54777520544776
$ws.Range("R2").Value = 805.92997684903
$ws.Range("S2").Value = 0.121746476107482
$ws.Range("T2").Value = 0.1309324101216233
$ws.Range("G3").Value = 2.454345666666666
$ws.Range("H3").Value = 7.363036999999999
$ws.Range("I3").Value = 0.2857469401257222
$ws.Range("J3").Value = 0.3027613480760606
$ws.Range("O3").Value = 0.1743777127077069
$ws.Range("P3").Value = 0.1769957741547643
$ws.Range("Q3").Value = 36.64974573976222
$ws.Range("R3").Value = 329.84771165786
$ws.Range("S3").Value = 0.04982789783234949
$ws.Range("T3").Value = 0.05358747918686241
$ws.Range("G4").Value = 2.454345666666666
$ws.Range("H4").Value = 7.363036999999999
$ws.Range("I4").Value = 0.2857469401257222
$ws.Range("J4").Value = 0.3027613480760606
$ws.Range("M4").Value = 9.680823666666667
$ws.Range("N4").Value = 29.042471
$ws.Range("O4").Value = 0.1130493445068016
$ws.Range("P4").Value = 0.1147466378470605
$ws.Range("Q4").Value = 23.76008761604744
$ws.Range("R4").Value = 213.840788544427
$ws.Range("S4").Value = 0.03230350427603718
$ws.Range("T4").Value = 0.03474084676177156
$ws.Range("G5").Value = 2.454345666666666
$ws.Range("H5").Value = 7.363036999999999
$ws.Range("I5").Value = 0.2857469401257222
$ws.Range("J5").Value = 0.3027613480760606
$ws.Range("M5").Value = 3.79999
$ws.Range("N5").Value = 7.59998
$ws.Range("O5").Value = 0.04437498227672168
$ws.Range("P5").Value = 0.0300274777826206
$ws.Range("Q5").Value = 9.326488989876665
$ws.Range("R5").Value = 55.95893393926
$ws.Range("S5").Value = 0.01268001540370637
$ws.Range("T5").Value = 0.009091159652790172
$ws.Range("G6").Value = 2.454345666666666
$ws.Range("H6").Value = 7.363036999999999
$ws.Range("I6").Value = 0.2857469401257222
$ws.Range("J6").Value = 0.3027613480760606
$ws.Range("M6").Value = 20.734808
$ws.Range("N6").Value = 62.204424
$ws.Range("O6").Value = 0.242133989171347
$ws.Range("P6").Value = 0.245769325661477
$ws.Range("Q6").Value = 50.89038616396533
$ws.Range("R6").Value = 458.013475475688
$ws.Range("S6").Value = 0.06918904650614716
$ws.Range("T6").Value = 0.07440945235301312
$ws.Range("G7").Value = 4.686805000000001
$ws.Range("I7").Value = 0.5456607868665887
$ws.Range("J7").Value = 0.5781514068052169
$ws.Range("M7").Value = 36.48539666666667
$ws.Range("N7").Value = 109.45619
$ws.Range("O7").Value = 0.4260639713374229
$ws.Range("P7").Value = 0.4324607845540777
$ws.Range("Q7").Value = 170.9999395243167
$ws.Range("R7").Value = 1538.99945571885
$ws.Range("S7").Value = 0.2324864018554819
$ws.Range("T7").Value = 0.2500278109780278
$ws.Range("G8").Value = 4.686805000000001
$ws.Range("I8").Value = 0.5456607868665887
$ws.Range("J8").Value = 0.5781514068052169
$ws.Range("O8").Value = 0.1743777127077069
$ws.Range("P8").Value = 0.1769957741547643
$ws.Range("Q8").Value = 69.98615309763335
$ws.Range("R8").Value = 629.8753778787001
$ws.Range("S8").Value = 0.09515107992808326
$ws.Range("T8").Value = 0.1023303558261554
$ws.Range("G9").Value = 4.686805000000001
$ws.Range("I9").Value = 0.5456607868665887
$ws.Range("J9").Value = 0.5781514068052169
$ws.Range("M9").Value = 9.680823666666667
$ws.Range("N9").Value = 29.042471
$ws.Range("O9").Value = 0.1130493445068016
$ws.Range("P9").Value = 0.1147466378470605
$ws.Range("Q9").Value = 45.37213276505167
$ws.Range("R9").Value = 408.349194885465
$ws.Range("S9").Value = 0.06168659427833344
$ws.Range("T9").Value = 0.06634093009744679
$ws.Range("G10").Value = 4.686805000000001
$ws.Range("I10").Value = 0.5456607868665887
$ws.Range("J10").Value = 0.5781514068052169
$ws.Range("M10").Value = 3.79999
$ws.Range("N10").Value = 7.59998
$ws.Range("O10").Value = 0.04437498227672168
$ws.Range("P10").Value = 0.0300274777826206
$ws.Range("Q10").Value = 17.80981213195
$ws.Range("R10").Value = 106.8588727917
$ws.Range("S10").Value = 0.02421368774630688
$ws.Range("T10").Value = 0.0173604285228345
$ws.Range("G11").Value = 4.686805000000001
$ws.Range("I11").Value = 0.5456607868665887
$ws.Range("J11").Value = 0.5781514068052169
$ws.Range("M11").Value = 20.734808
$ws.Range("N11").Value = 62.204424
$ws.Range("O11").Value = 0.242133989171347
$ws.Range("P11").Value = 0.245769325661477
$ws.Range("Q11").Value = 97.18000180844001
$ws.Range("R11").Value = 874.6200162759601
$ws.Range("S11").Value = 0.1321230230583833
$ws.Range("T11").Value = 0.1420918813807524
$ws.Range("G12").Value = 1.4480775
$ws.Range("H12").Value = 2.896155
$ws.Range("I12").Value = 0.1685922730076891
$ws.Range("J12").Value = 0.1190872451187225
$ws.Range("M12").Value = 36.48539666666667
$ws.Range("N12").Value = 109.45619
$ws.Range("O12").Value = 0.4260639713374229
$ws.Range("P12").Value = 0.4324607845540777
$ws.Range("Q12").Value = 52.833681991575
$ws.Range("R12").Value = 317.00209194945
$ws.Range("S12").Value = 0.07183109337445903
$ws.Range("T12").Value = 0.05150056345442649
$ws.Range("G13").Value = 1.4480775
$ws.Range("H13").Value = 2.896155
$ws.Range("I13").Value = 0.1685922730076891
$ws.Range("J13").Value = 0.1190872451187225
$ws.Range("O13").Value = 0.1743777127077069
$ws.Range("P13").Value = 0.1769957741547643
$ws.Range("Q13").Value = 21.62355242265
$ws.Range("R13").Value = 129.7413145359
$ws.Range("S13").Value = 0.0293987349472741
$ws.Range("T13").Value = 0.02107793914174647
$ws.Range("G14").Value = 1.4480775
$ws.Range("H14").Value = 2.896155
$ws.Range("I14").Value = 0.1685922730076891
$ws.Range("J14").Value = 0.1190872451187225
$ws.Range("M14").Value = 9.680823666666667
$ws.Range("N14").Value = 29.042471
$ws.Range("O14").Value = 0.1130493445068016
$ws.Range("P14").Value = 0.1147466378470605
$ws.Range("Q14").Value = 14.0185829331675
$ws.Range("R14").Value = 84.11149759900499
$ws.Range("S14").Value = 0.019059245952431
$ws.Range("T14").Value = 0.01366486098784218
$ws.Range("G15").Value = 1.4480775
$ws.Range("H15").Value = 2.896155
$ws.Range("I15").Value = 0.1685922730076891
$ws.Range("J15").Value = 0.1190872451187225
$ws.Range("M15").Value = 3.79999
$ws.Range("N15").Value = 7.59998
$ws.Range("O15").Value = 0.04437498227672168
$ws.Range("P15").Value = 0.0300274777826206
$ws.Range("Q15").Value = 5.502680019225
$ws.Range("R15").Value = 22.0107200769
$ws.Range("S15").Value = 0.007481279126708429
$ws.Range("T15").Value = 0.003575889606995933
$ws.Range("G16").Value = 1.4480775
$ws.Range("H16").Value = 2.896155
$ws.Range("I16").Value = 0.1685922730076891
$ws.Range("J16").Value = 0.1190872451187225
$ws.Range("M16").Value = 20.734808
$ws.Range("N16").Value = 62.204424
$ws.Range("O16").Value = 0.242133989171347
$ws.Range("P16").Value = 0.245769325661477
$ws.Range("Q16").Value = 30.02560893162
$ws.Range("R16").Value = 180.15365358972
$ws.Range("S16").Value = 0.04082191960681658
$ws.Range("T16").Value = 0.02926799192771144
